$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '99.290.50'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '3.309.34'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.92'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '628.12'
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  +21.36%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.411'
$ws.Range("E8").Value = '  +6.88%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.994'
$ws.Range("E10").Value = '  +23.72%  '
$ws.Range("D11").Value = '3.310.53'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("E12").Value = '  +2.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.33'
$ws.Range("E13").Value = '  +15.24%  '
$ws.Range("D14").Value = '99.005.31'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000253'
$ws.Range("E15").Value = '  +2.39%  '
$ws.Range("D16").Value = '3.932.38'
$ws.Range("E16").Value = '  -1.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.49'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").Value = '3.308.59'
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.49'
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.74'
$ws.Range("E20").Value = '  +6.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.45'
$ws.Range("E21").Value = '  +9.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '487.07'
$ws.Range("E22").Value = '  +1.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.51'
$ws.Range("E23").Value = '  +3.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000204'
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.79'
$ws.Range("E25").Value = '  +1.71%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.340'
$ws.Range("E26").Value = '  +35.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.05'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.24'
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("D29").Value = '3.485.04'
$ws.Range("E29").Value = '  -2.41%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.154'
$ws.Range("E30").Value = '  +24.15%  '
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.192'
$ws.Range("E32").Value = '  +1.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.66'
$ws.Range("E33").Value = '  +15.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.05'
$ws.Range("E35").Value = '  +3.15%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.153'
$ws.Range("E36").Value = '  +1.20%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.484'
$ws.Range("E37").Value = '  +7.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.40'
$ws.Range("E38").Value = '  +1.79%  '
$ws.Range("E39").Value = '  +1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '499.43'
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.75'
$ws.Range("E41").Value = '  -0.24%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.90'
$ws.Range("E42").Value = '  +6.31%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.789'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.23'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '160.28'
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.96'
$ws.Range("E48").Value = '  +2.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.83'
$ws.Range("E49").Value = '  +6.95%  '
$ws.Range("E50").Value = '  +7.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.35'
$ws.Range("E51").Value = '  +15.04%  '
